# Added DAP file for today
# Updates TOTAL_SS_LOAD (column S) values for hours 10-22 (rows 11-23),
# and refreshes the WESM_RATE / CURRENT_RATE (columns V/W) spot that
# used to sit on row 12 but now belongs on row 23 with today's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TOTAL_SS_LOAD (column S) updates
$ws.Range("S11").Value = 32621
$ws.Range("S12").Value = 33187
$ws.Range("S13").Value = 32997
$ws.Range("S14").Value = 34647
$ws.Range("S15").Value = 35793
$ws.Range("S16").Value = 37060
$ws.Range("S17").Value = 37376
$ws.Range("S18").Value = 37480
$ws.Range("S19").Value = 40256
$ws.Range("S20").Value = 42234
$ws.Range("S21").Value = 40306
$ws.Range("S22").Value = 40825
$ws.Range("S23").Value = 41147

# WESM_RATE / CURRENT_RATE (columns V/W) move from row 12 to row 23
$ws.Range("V12").ClearContents()
$ws.Range("W12").ClearContents()

$ws.Range("V23").Value = 3927.092180555555
$ws.Range("W23").Value = 7.757045406664108
